$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 15770.80323680944
$ws.Range("D3").Value = 1052.413686175801

$ws.Range("B4").Value = 6632.685490293635
$ws.Range("D4").Value = 798.7814462204649

$ws.Range("B5").Value = 4768.756578082192
$ws.Range("D5").Value = 37.26064657534252

$ws.Range("B6").Value = 10357.34987534247
$ws.Range("D6").Value = 265.4806212328762

$ws.Range("B7").Value = 13674.28660273972
$ws.Range("D7").Value = 959.2604273972595

$ws.Range("B8").Value = 21235.57258630153
$ws.Range("D8").Value = 1920

$ws.Range("B9").Value = 29926.40090684943
$ws.Range("D9").Value = 1920.000060273973

$ws.Range("F10").Value = 22709956.66100555

$ws.Range("G11").Value = 0.7243161176691555

$ws.Range("F12").Value = 1515475.708093151
$ws.Range("G12").Value = 0.06673177455663397

$ws.Range("G13").Value = 0.2089521077742105
